$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove D128 (no longer present after the edit)
$ws.Cells.Item(128, 4).ClearContents()

# Update existing numeric values (re-computed ifoCAST error series)
$ws.Cells.Item(129, 4).Value = 0.7804486234241067
$ws.Cells.Item(130, 4).Value = 0.7918746594241067
$ws.Cells.Item(131, 4).Value = 0.7115302104241067
$ws.Cells.Item(132, 4).Value = 0.7732527034241068
$ws.Cells.Item(133, 3).Value = 0.5947585844241068
$ws.Cells.Item(134, 3).Value = -0.1471494035758933
$ws.Cells.Item(135, 3).Value = 0.2315426864241067
$ws.Cells.Item(136, 3).Value = 0.5186180304241067
$ws.Cells.Item(137, 3).Value = 0.5462623554241067
$ws.Cells.Item(138, 3).Value = 0.3054124294241067
$ws.Cells.Item(139, 3).Value = 0.3912781354241067
$ws.Cells.Item(140, 2).Value = 0.0999529544241067

# Append new sampling rows (ifoCAST-sampling), copying the date-label
# formatting used throughout column A
$ws.Cells.Item(140, 1).Copy()
$ws.Cells.Item(141, 1).PasteSpecial(-4122)
$ws.Cells.Item(142, 1).PasteSpecial(-4122)
$ws.Cells.Item(143, 1).PasteSpecial(-4122)
$ws.Cells.Item(144, 1).PasteSpecial(-4122)
$ws.Cells.Item(145, 1).PasteSpecial(-4122)

$ws.Cells.Item(141, 1).Value = "2025-07-25_diff"
$ws.Cells.Item(142, 1).Value = "2025-08-07_diff"
$ws.Cells.Item(143, 1).Value = "2025-08-22_diff"
$ws.Cells.Item(144, 1).Value = "2025-08-25_diff"
$ws.Cells.Item(145, 1).Value = "2025-09-08_diff"

$ws.Cells.Item(141, 2).Value = 0.05603945542410671
